$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values like "0.9997" or
# "1.000" are not auto-converted to numbers by the smart input parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.000.61'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.901.05'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '0.7879'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '244.56'
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.3161'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '25.87'
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").Value = '0.07319'
$ws.Range("E10").Value = '  +4.57%  '
$ws.Range("D11").Value = '0.08122'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").Value = '0.7792'
$ws.Range("E12").Value = '  +2.88%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.514'
$ws.Range("E13").Value = '  +4.24%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.854.07'
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '94.35'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '6.263'
$ws.Range("E16").Value = '  +5.81%  '
$ws.Range("D17").Value = '29.964.23'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '14.02'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("D19").Value = '246.96'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '0.000007846'
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '8.181'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '2.115.31'
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '0.1607'
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").Value = '9.498'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("D27").Value = '163.34'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("D28").Value = '18.83'
$ws.Range("E28").Value = '  +1.23%  '
$ws.Range("D29").Value = '2.045'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").Value = '1.441'
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").Value = '4.499'
$ws.Range("E32").Value = '  +2.77%  '
$ws.Range("D33").Value = '0.05620'
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").Value = '4.104'
$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("D35").Value = '1.254'
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").Value = '0.7564'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("D37").Value = '1.002'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = '2.677'
$ws.Range("E38").Value = '  +3.24%  '
$ws.Range("D39").Value = '0.01939'
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("D40").Value = '2.796'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").Value = '1.148.58'
$ws.Range("E41").Value = '  +12.83%  '
$ws.Range("D42").Value = '0.4483'
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("D43").Value = '74.10'
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("D44").Value = '5.990'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("D45").Value = '0.8571'
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("D46").Value = '1.908'
$ws.Range("E46").Value = '  +2.69%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '0.9998'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").Value = '3.166'
$ws.Range("E48").Value = '  +9.03%  '
$ws.Range("D49").Value = '102.35'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.827'
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '7.555'
$ws.Range("E51").Value = '  +1.72%  '

# Restore the original (default) cell style now that the text values are set,
# so no lingering text-formatting style is left on the Price column.
$ws.Range("D2:D51").Style = "Normal"
